{"js": "// AT01 Q5 Code of Ethics \u2014 expand the \"Introduction: \" paragraph into the\n// full introduction text, written as a sequence of runs (mirroring how Word\n// itself breaks a freshly-typed passage into runs around autocorrect /\n// grammar-checker marks such as the \"Co\" abbreviation flagged by proofing).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the paragraph that currently holds \"Introduction: \" (style \"My Style\").\nlet introParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Introduction: \") {\n    introParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!introParagraph) {\n  throw new Error('Could not find the \"Introduction: \" paragraph.');\n}\n\nconst introRange = introParagraph.getRange();\n\n// Read back the paragraph's own OOXML so we can keep its identity\n// (w14:paraId/rsid*) and paragraph properties (pStyle, rPr) untouched \u2014\n// only the runs inside it are being replaced.\nconst ooxml = introRange.getOoxml();\nawait context.sync();\nconst fullOoxml = ooxml.value;\n\nconst pOpenTagMatch = fullOoxml.match(/<w:p(?:\\s[^>]*)?>/);\nconst pOpenTag = pOpenTagMatch ? pOpenTagMatch[0] : \"<w:p>\";\n\nconst pPrMatch = fullOoxml.match(/<w:pPr>[\\s\\S]*?<\\/w:pPr>/);\nconst pPr = pPrMatch ? pPrMatch[0] : \"\";\n\n// The new text, split the way the live edit produced it: a run per typing/\n// autocorrect chunk, with the abbreviation \"Co\" wrapped in proofErr marks\n// the way Word's grammar checker flags it (\"gramStart\"/\"gramEnd\").\nconst segments = [\n  { text: \"Here at\" },\n  { text: \" Rainbow Hero \" },\n  { text: \"Co\", gramFlag: true },\n  { text: \" we uphold\" },\n  { text: \" a high level\" },\n  { text: \" \" },\n  { text: \"of \" },\n  { text: \"Honesty\" },\n  { text: \". We \" },\n  { text: \"prioritize\" },\n  { text: \" on\" },\n  { text: \" \" },\n  { text: \"Privacy, \" },\n  { text: \"Copyright, Health & Safety in the \" },\n  { text: \"workplace. At Rainbow Hero Co we enjoy have diversity & want to make sure everyone can be included in the tasks. \" },\n  { text: \"Our corporation encourages everyone to get more Professional Development when possible.\" },\n];\n\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\");\n}\n\nlet runsXml = \"\";\nsegments.forEach((seg, idx) => {\n  const needsPreserve = seg.text !== seg.text.trim() || seg.text.length === 0;\n  const preserveAttr = needsPreserve ? ' xml:space=\"preserve\"' : \"\";\n  // The very first run reuses the original run's rsidRPr; the rest are\n  // freshly minted plain runs, same as they came from separate edits.\n  const openRunTag = idx === 0 ? '<w:r w:rsidRPr=\"00D659BD\">' : \"<w:r>\";\n  const runXml =\n    openRunTag +\n    '<w:rPr><w:lang w:val=\"en-US\"/></w:rPr>' +\n    `<w:t${preserveAttr}>${xmlEscape(seg.text)}</w:t>` +\n    \"</w:r>\";\n  runsXml += seg.gramFlag\n    ? '<w:proofErr w:type=\"gramStart\"/>' + runXml + '<w:proofErr w:type=\"gramEnd\"/>'\n    : runXml;\n});\n\nconst newParagraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  \"<w:body>\" +\n  pOpenTag +\n  pPr +\n  runsXml +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nintroRange.insertOoxml(newParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# AT01 Q5 Code of Ethics -- expand the \"Introduction: \" paragraph into the\n# full introduction text, written as a sequence of runs (mirroring how Word\n# itself breaks a freshly-typed passage into runs around autocorrect /\n# grammar-checker marks such as the \"Co\" abbreviation flagged by proofing).\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that currently holds \"Introduction: \" (style \"My Style\").\n$target = $null\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n  $candidate = $paras.Item($i)\n  if ($candidate.Range.Text.StartsWith(\"Introduction:\")) {\n    $target = $candidate\n    break\n  }\n}\nif ($target -eq $null) {\n  throw \"Could not find the 'Introduction: ' paragraph.\"\n}\n\n$r = $target.Range\n\n# Read back the paragraph's own OOXML so we can keep its identity\n# (w14:paraId/rsid*) and paragraph properties (pStyle, rPr) untouched --\n# only the runs inside it are being replaced.\n$xml = $r.WordOpenXML\n\n$pOpenTag = \"<w:p>\"\nif ($xml -match '<w:p(?: [^>]*)?>') {\n  $pOpenTag = $matches[0]\n}\n$pPr = \"\"\nif ($xml -match '<w:pPr>[\\s\\S]*?</w:pPr>') {\n  $pPr = $matches[0]\n}\n\n# The new text, split the way the live edit produced it: a run per typing/\n# autocorrect chunk, with the abbreviation \"Co\" wrapped in proofErr marks\n# the way Word's grammar checker flags it (\"gramStart\"/\"gramEnd\").\n$runsXml = '<w:r w:rsidRPr=\"00D659BD\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Here at</w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> Rainbow Hero </w:t></w:r>'\n$runsXml += '<w:proofErr w:type=\"gramStart\"/>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Co</w:t></w:r>'\n$runsXml += '<w:proofErr w:type=\"gramEnd\"/>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> we uphold</w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> a high level</w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">of </w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Honesty</w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">. We </w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>prioritize</w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> on</w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">Privacy, </w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">Copyright, Health &amp; Safety in the </w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">workplace. At Rainbow Hero Co we enjoy have diversity &amp; want to make sure everyone can be included in the tasks. </w:t></w:r>'\n$runsXml += '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Our corporation encourages everyone to get more Professional Development when possible.</w:t></w:r>'\n\n$newParagraphOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>'\n$newParagraphOoxml += '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">'\n$newParagraphOoxml += '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">'\n$newParagraphOoxml += '<pkg:xmlData>'\n$newParagraphOoxml += '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">'\n$newParagraphOoxml += '<w:body>'\n$newParagraphOoxml += $pOpenTag\n$newParagraphOoxml += $pPr\n$newParagraphOoxml += $runsXml\n$newParagraphOoxml += '</w:p>'\n$newParagraphOoxml += '</w:body>'\n$newParagraphOoxml += '</w:document>'\n$newParagraphOoxml += '</pkg:xmlData>'\n$newParagraphOoxml += '</pkg:part>'\n$newParagraphOoxml += '</pkg:package>'\n\n$r.InsertXML($newParagraphOoxml)\n"}
